$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B7, F7, B12 from "special" (30) back to "time" (10)
$ws.Range("B7").Value = "time"
$ws.Range("F7").Value = "time"
$ws.Range("B12").Value = "time"

# Change the selection to A3:B3
$ws.Range("A3:B3").Select()
